# Insert a new "Jurisdiction" property row into the Metadata sheet,
# directly after the existing "Contact" row (row 10) and before
# "Description" (previously row 11), pushing Description/Purpose/
# Copyright/Immutable down by one row (11->12, 12->13, 13->14, 14->15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Snapshot the existing values for rows 11-14 before they shift down,
# since each cell's .Value must be read as a method call to get the
# actual stored data (not the property descriptor).
$a11 = $ws.Cells.Item(11, 1).Value()
$b11 = $ws.Cells.Item(11, 2).Value()
$a12 = $ws.Cells.Item(12, 1).Value()
$b12 = $ws.Cells.Item(12, 2).Value()
$a13 = $ws.Cells.Item(13, 1).Value()
$b13 = $ws.Cells.Item(13, 2).Value()
$a14 = $ws.Cells.Item(14, 1).Value()
$b14 = $ws.Cells.Item(14, 2).Value()

# Row 15 is brand new - clone the formatting used by the row above it
# (the common body-row style) before writing the shifted-down values.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 down into 12-15.
$ws.Cells.Item(15, 1).Value = $a14
$ws.Cells.Item(15, 2).Value = $b14
$ws.Cells.Item(14, 1).Value = $a13
$ws.Cells.Item(14, 2).Value = $b13
$ws.Cells.Item(13, 1).Value = $a12
$ws.Cells.Item(13, 2).Value = $b12
$ws.Cells.Item(12, 1).Value = $a11
$ws.Cells.Item(12, 2).Value = $b11

# Write the new Jurisdiction property into the now-vacated row 11; it
# keeps the body-row style that row already carried.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
